# Deploy the implementation guide.
# Regenerated CodeSystem metadata sheet: refreshed the build Date, populated
# the publisher Contact display text, and added a new Jurisdiction property
# row (FHIR IG Publisher now emits this row even though no jurisdictions are
# configured, hence the empty value).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# 1) Refresh the build Date value.
$ws1.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# 2) Contact now renders with a publisher name + URL instead of the
#    "no display" placeholder.
$ws1.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# 3) Insert a new "Jurisdiction" row directly below "Contact" (row 11),
#    pushing "Description" and everything after it down by one row.
$ws1.Rows.Item(11).Insert()

# Populate the new row's cells before re-applying formatting so the cell
# holding an explicit empty string keeps a plain (non quote-prefixed) style,
# matching the other data rows.
$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = "'"

# Re-apply the same formatting (fill/border/alignment) used by the other
# metadata rows to the freshly inserted row, since a bare row Insert() does
# not inherit it.
$ws1.Range("A12:B12").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
